# LogDolar.xlsx update: new "Dia- Hora Ultima Actualizacion" header text,
# re-style the two existing closed-market log rows (38-39) to match the
# rest of the log table, and append five new log rows (40-44) capturing
# the latest real-time dollar quotes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header -----------------------------------------------------------
$ws.Range("A1").Value = "Dia- Hora Ultima Actualizacion "

# --- Re-format rows 38-39 (already-logged rows) so they match the rest
# of the table's look (border/fill styling carried by rows 2-37). -----
$ws.Range("A38:C39").Style = $ws.Range("A37:C37").Style
$ws.Range("D38:D39").Style = $ws.Range("D37").Style

# --- New log rows -------------------------------------------------------
$ws.Range("A40").Value = "12:00:09"
$ws.Range("B40").Value = "Info en tiempo real"
$ws.Range("C40").Value = "02/01/2023 04:49"
$ws.Range("D40").Value = "848,25"

$ws.Range("A41").Value = "12:00:09"
$ws.Range("B41").Value = "Info en tiempo real"
$ws.Range("C41").Value = "02/01/2023 04:51"
$ws.Range("D41").Value = "848,25"

$ws.Range("A42").Value = "12:00:09"
$ws.Range("B42").Value = "Info en tiempo real"
$ws.Range("C42").Value = "02/01/2023 04:54"
$ws.Range("D42").Value = "848,25"

$ws.Range("A40:D42").Style = $ws.Range("A37:D37").Style

$ws.Range("A43").Value = "23:00:03"
$ws.Range("B43").Value = "Info en tiempo real"
$ws.Range("C43").Value = "02/01/2023 05:14"
$ws.Range("D43").Value = "848,25"

$ws.Range("A44").Value = "23:00:03"
$ws.Range("B44").Value = "Info en tiempo real"
$ws.Range("C44").Value = "02/01/2023 05:16"
$ws.Range("D44").Value = "848,25"

# --- Cosmetic: restore the small page-margin tweak seen after reopening
# the file (header/footer margin in points; XML stores inches). --------
$ws.PageSetup.HeaderMargin = 36.850393700787386
$ws.PageSetup.FooterMargin = 36.850393700787386

# --- Selection as left by the editing session --------------------------
$ws.Range("B7").Select()
